$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 457.63635
$ws.Range("I15").Value = 457.63635
$ws.Range("K15").Value = 1372.90905
$ws.Range("M15").Value = -1203.90905

$ws.Range("H40").Value = 2043.5652
$ws.Range("J40").Value = 2050.1
$ws.Range("L40").Value = 2050.1
$ws.Range("N40").Value = -2400.1

$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 6000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5516

$ws.Range("H61").Value = 7657
$ws.Range("I61").Value = 899.75
$ws.Range("J61").Value = 16666.666
$ws.Range("K61").Value = 2699.25
$ws.Range("L61").Value = 49999.99800000001
$ws.Range("M61").Value = -2527.25
$ws.Range("N61").Value = -50343.99800000001

$ws.Range("H94").Value = 879.75
$ws.Range("I94").Value = 996.55554
$ws.Range("K94").Value = 996.55554
$ws.Range("M94").Value = -545.55554

$ws.Range("H132").Value = 1066.8269
$ws.Range("I132").Value = 982.5106
$ws.Range("K132").Value = 2947.5318
$ws.Range("M132").Value = -417.5317999999997

$ws.Range("H137").Value = 9675.352999999999
$ws.Range("I137").Value = 4296.564
$ws.Range("K137").Value = 12889.692
$ws.Range("M137").Value = -10339.692

$ws.Range("H138").Value = 4701.8237
$ws.Range("I138").Value = 3515.074
$ws.Range("J138").Value = 5483.3413
$ws.Range("K138").Value = 10545.222
$ws.Range("L138").Value = 16450.0239
$ws.Range("M138").Value = -5405.222
$ws.Range("N138").Value = -26730.0239

$ws.Range("H141").Value = 1876
$ws.Range("I141").Value = 1953.1
$ws.Range("J141").Value = 1105
$ws.Range("K141").Value = 5859.299999999999
$ws.Range("L141").Value = 3315
$ws.Range("M141").Value = -679.2999999999993
$ws.Range("N141").Value = -13675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3887.1714
$ws.Range("I32").Value = 2339.7354
$ws.Range("K32").Value = 2339.7354
$ws.Range("M32").Value = -2052.7354

$ws.Range("H74").Value = 5009.8
$ws.Range("I74").Value = 3372.3333
$ws.Range("J74").Value = 7466
$ws.Range("K74").Value = 3372.3333
$ws.Range("L74").Value = 7466
$ws.Range("M74").Value = -2498.3333
$ws.Range("N74").Value = -9214

$ws.Range("H77").Value = 5009.8
$ws.Range("I77").Value = 3372.3333
$ws.Range("J77").Value = 7466
$ws.Range("K77").Value = 16861.6665
$ws.Range("L77").Value = 37330
$ws.Range("M77").Value = -12493.6665
$ws.Range("N77").Value = -46066

$ws.Range("H132").Value = 5237.1665
$ws.Range("I132").Value = 3571.2104
$ws.Range("J132").Value = 11567.8
$ws.Range("K132").Value = 10713.6312
$ws.Range("L132").Value = 34703.39999999999
$ws.Range("M132").Value = -8183.6312
$ws.Range("N132").Value = -39763.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2639.25
$ws.Range("I99").Value = 2074.1304
$ws.Range("K99").Value = 2074.1304
$ws.Range("M99").Value = -576.1304

$ws.Range("H105").Value = 4487.1665
$ws.Range("I105").Value = 4981
$ws.Range("K105").Value = 4981
$ws.Range("M105").Value = -3234

$ws.Range("H107").Value = 860.7
$ws.Range("I107").Value = 908.1429000000001
$ws.Range("K107").Value = 908.1429000000001
$ws.Range("M107").Value = 1011.8571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 25258.842
$ws.Range("I58").Value = 30719.059
$ws.Range("J58").Value = 6694.1
$ws.Range("K58").Value = 30719.059
$ws.Range("L58").Value = 6694.1
$ws.Range("M58").Value = -30516.059
$ws.Range("N58").Value = -7100.1

$ws.Range("H59").Value = 33459.266
$ws.Range("J59").Value = 30530.076
$ws.Range("L59").Value = 30530.076
$ws.Range("N59").Value = -32820.076

$ws.Range("H62").Value = 341496.56
$ws.Range("J62").Value = 295138.56
$ws.Range("L62").Value = 295138.56
$ws.Range("N62").Value = -296386.56

$ws.Range("H65").Value = 341496.56
$ws.Range("J65").Value = 295138.56
$ws.Range("L65").Value = 1475692.8
$ws.Range("N65").Value = -1481932.8

$ws.Range("H122").Value = 1055.1904
$ws.Range("I122").Value = 965.41174
$ws.Range("J122").Value = 1436.75
$ws.Range("K122").Value = 2896.23522
$ws.Range("L122").Value = 4310.25
$ws.Range("M122").Value = -446.23522
$ws.Range("N122").Value = -9210.25

$ws.Range("H132").Value = 19526.04
$ws.Range("I132").Value = 12364.019
$ws.Range("K132").Value = 37092.057
$ws.Range("M132").Value = -34562.057

$ws.Range("H134").Value = 3992.82
$ws.Range("I134").Value = 2025.4147
$ws.Range("J134").Value = 12955.444
$ws.Range("K134").Value = 6076.2441
$ws.Range("L134").Value = 38866.33199999999
$ws.Range("M134").Value = -3541.2441
$ws.Range("N134").Value = -43936.33199999999

$ws.Range("H136").Value = 25258.842
$ws.Range("I136").Value = 30719.059
$ws.Range("J136").Value = 6694.1
$ws.Range("K136").Value = 92157.177
$ws.Range("L136").Value = 20082.3
$ws.Range("M136").Value = -89607.177
$ws.Range("N136").Value = -25182.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1149.6046
$ws.Range("I5").Value = 503
$ws.Range("K5").Value = 1509
$ws.Range("M5").Value = -1397

$ws.Range("H60").Value = 2749
$ws.Range("I60").Value = 780.3333
$ws.Range("K60").Value = 2340.9999
$ws.Range("M60").Value = -2089.9999

$ws.Range("H113").Value = 1364.5
$ws.Range("J113").Value = 1130.875
$ws.Range("L113").Value = 3392.625
$ws.Range("N113").Value = -7732.625

$ws.Range("H121").Value = 1946.1034
$ws.Range("I121").Value = 1410.5217
$ws.Range("K121").Value = 4231.5651
$ws.Range("M121").Value = -2921.5651

$ws.Range("H122").Value = 16666950
$ws.Range("J122").Value = 20000160
$ws.Range("L122").Value = 180001440
$ws.Range("N122").Value = -180006340

$ws.Range("H131").Value = 1904.931
$ws.Range("J131").Value = 2062.353
$ws.Range("L131").Value = 6187.059
$ws.Range("N131").Value = -16267.059

$ws.Range("H135").Value = 1149.6046
$ws.Range("I135").Value = 503
$ws.Range("K135").Value = 4527
$ws.Range("M135").Value = -1992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 751.12
$ws.Range("I97").Value = 813.2857
$ws.Range("K97").Value = 813.2857
$ws.Range("M97").Value = -317.2857

$ws.Range("H126").Value = 3594
$ws.Range("I126").Value = 3594
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10782
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8312
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1475
$ws.Range("I13").Value = 1475
$ws.Range("K13").Value = 1475
$ws.Range("M13").Value = -1335

$ws.Range("H40").Value = 1715.08
$ws.Range("I40").Value = 1453.0769
$ws.Range("K40").Value = 1453.0769
$ws.Range("M40").Value = -1317.0769

$ws.Range("H46").Value = 1500.1
$ws.Range("I46").Value = 1066.6666
$ws.Range("J46").Value = 1608.4584
$ws.Range("K46").Value = 1066.6666
$ws.Range("L46").Value = 1608.4584
$ws.Range("M46").Value = -878.6666
$ws.Range("N46").Value = -1984.4584

$ws.Range("H61").Value = 1804.6
$ws.Range("I61").Value = 1659.6154
$ws.Range("K61").Value = 1659.6154
$ws.Range("M61").Value = -1457.6154

$ws.Range("H113").Value = 1804.6
$ws.Range("I113").Value = 1659.6154
$ws.Range("K113").Value = 1659.6154
$ws.Range("M113").Value = 510.3846000000001

$ws.Range("H122").Value = 2429.6667
$ws.Range("I122").Value = 2472.389
$ws.Range("K122").Value = 7417.167
$ws.Range("M122").Value = -4967.167

$ws.Range("H132").Value = 5546.5386
$ws.Range("I132").Value = 5748.6665
$ws.Range("J132").Value = 5091.75
$ws.Range("K132").Value = 17245.9995
$ws.Range("L132").Value = 15275.25
$ws.Range("M132").Value = -14715.9995
$ws.Range("N132").Value = -20335.25

$ws.Range("H134").Value = 29999
$ws.Range("J134").Value = 29999
$ws.Range("L134").Value = 29999
$ws.Range("N134").Value = -40139

$ws.Range("H136").Value = 2754.878
$ws.Range("I136").Value = 2248.2856
$ws.Range("J136").Value = 3846
$ws.Range("K136").Value = 6744.8568
$ws.Range("L136").Value = 11538
$ws.Range("M136").Value = -4194.8568
$ws.Range("N136").Value = -16638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4314.9165
$ws.Range("I126").Value = 5165.931
$ws.Range("J126").Value = 789.2857
$ws.Range("K126").Value = 15497.793
$ws.Range("L126").Value = 2367.8571
$ws.Range("M126").Value = -13027.793
$ws.Range("N126").Value = -7307.8571

$ws.Range("H132").Value = 21279.865
$ws.Range("I132").Value = 11683.147
$ws.Range("J132").Value = 39407
$ws.Range("K132").Value = 35049.44100000001
$ws.Range("L132").Value = 118221
$ws.Range("M132").Value = -32519.44100000001
$ws.Range("N132").Value = -123281
